$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lgi4"
$ws.Cells.Item(2, 3).Value = "Adam22"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1241926666666667
$ws.Cells.Item(2, 8).Value = 0.372578
$ws.Cells.Item(2, 9).Value = 0.00757902233016378
$ws.Cells.Item(2, 10).Value = 0.00757902233016378
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8328803333333332
$ws.Cells.Item(2, 14).Value = 2.498641
$ws.Cells.Item(2, 15).Value = 0.03917234793046614
$ws.Cells.Item(2, 16).Value = 0.03917234793046613
$ws.Cells.Item(2, 17).Value = 0.1034376296108889
$ws.Cells.Item(2, 18).Value = 0.9309386664979998
$ws.Cells.Item(2, 19).Value = 0.0002968880996899478
$ws.Cells.Item(2, 20).Value = 0.0002968880996899477

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lgi4"
$ws.Cells.Item(3, 3).Value = "Adam22"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1241926666666667
$ws.Cells.Item(3, 8).Value = 0.372578
$ws.Cells.Item(3, 9).Value = 0.00757902233016378
$ws.Cells.Item(3, 10).Value = 0.00757902233016378
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 4.873409
$ws.Cells.Item(3, 14).Value = 14.620227
$ws.Cells.Item(3, 15).Value = 0.2292080450398417
$ws.Cells.Item(3, 16).Value = 0.2292080450398417
$ws.Cells.Item(3, 17).Value = 0.6052416594673332
$ws.Cells.Item(3, 18).Value = 5.447174935205999
$ws.Cells.Item(3, 19).Value = 0.001737172891610146
$ws.Cells.Item(3, 20).Value = 0.001737172891610146

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lgi4"
$ws.Cells.Item(4, 3).Value = "Adam22"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1241926666666667
$ws.Cells.Item(4, 8).Value = 0.372578
$ws.Cells.Item(4, 9).Value = 0.00757902233016378
$ws.Cells.Item(4, 10).Value = 0.00757902233016378
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 15.55565633333333
$ws.Cells.Item(4, 14).Value = 46.66696899999999
$ws.Cells.Item(4, 15).Value = 0.7316196070296923
$ws.Cells.Item(4, 16).Value = 0.7316196070296922
$ws.Cells.Item(4, 17).Value = 1.931898441786889
$ws.Cells.Item(4, 18).Value = 17.38708597608199
$ws.Cells.Item(4, 19).Value = 0.005544961338863688
$ws.Cells.Item(4, 20).Value = 0.005544961338863687

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lgi4"
$ws.Cells.Item(5, 3).Value = "Adam22"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.566363666666668
$ws.Cells.Item(5, 8).Value = 25.699091
$ws.Cells.Item(5, 9).Value = 0.5227737133000635
$ws.Cells.Item(5, 10).Value = 0.5227737133000636
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8328803333333332
$ws.Cells.Item(5, 14).Value = 2.498641
$ws.Cells.Item(5, 15).Value = 0.03917234793046614
$ws.Cells.Item(5, 16).Value = 0.03917234793046613
$ws.Cells.Item(5, 17).Value = 7.134755826147889
$ws.Cells.Item(5, 18).Value = 64.212802435331
$ws.Cells.Item(5, 19).Value = 0.02047827378629184
$ws.Cells.Item(5, 20).Value = 0.02047827378629184

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lgi4"
$ws.Cells.Item(6, 3).Value = "Adam22"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 8.566363666666668
$ws.Cells.Item(6, 8).Value = 25.699091
$ws.Cells.Item(6, 9).Value = 0.5227737133000635
$ws.Cells.Item(6, 10).Value = 0.5227737133000636
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 4.873409
$ws.Cells.Item(6, 14).Value = 14.620227
$ws.Cells.Item(6, 15).Value = 0.2292080450398417
$ws.Cells.Item(6, 16).Value = 0.2292080450398417
$ws.Cells.Item(6, 17).Value = 41.74739379040633
$ws.Cells.Item(6, 18).Value = 375.726544113657
$ws.Cells.Item(6, 19).Value = 0.1198239408237263
$ws.Cells.Item(6, 20).Value = 0.1198239408237263

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lgi4"
$ws.Cells.Item(7, 3).Value = "Adam22"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 8.566363666666668
$ws.Cells.Item(7, 8).Value = 25.699091
$ws.Cells.Item(7, 9).Value = 0.5227737133000635
$ws.Cells.Item(7, 10).Value = 0.5227737133000636
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 15.55565633333333
$ws.Cells.Item(7, 14).Value = 46.66696899999999
$ws.Cells.Item(7, 15).Value = 0.7316196070296923
$ws.Cells.Item(7, 16).Value = 0.7316196070296922
$ws.Cells.Item(7, 17).Value = 133.2554092250199
$ws.Cells.Item(7, 18).Value = 1199.298683025179
$ws.Cells.Item(7, 19).Value = 0.3824714986900455
$ws.Cells.Item(7, 20).Value = 0.3824714986900455

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lgi4"
$ws.Cells.Item(8, 3).Value = "Adam22"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 7.695813999999999
$ws.Cells.Item(8, 8).Value = 23.087442
$ws.Cells.Item(8, 9).Value = 0.4696472643697726
$ws.Cells.Item(8, 10).Value = 0.4696472643697726
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.8328803333333332
$ws.Cells.Item(8, 14).Value = 2.498641
$ws.Cells.Item(8, 15).Value = 0.03917234793046614
$ws.Cells.Item(8, 16).Value = 0.03917234793046613
$ws.Cells.Item(8, 17).Value = 6.409692129591332
$ws.Cells.Item(8, 18).Value = 57.68722916632199
$ws.Cells.Item(8, 19).Value = 0.01839718604448435
$ws.Cells.Item(8, 20).Value = 0.01839718604448434

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lgi4"
$ws.Cells.Item(9, 3).Value = "Adam22"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 7.695813999999999
$ws.Cells.Item(9, 8).Value = 23.087442
$ws.Cells.Item(9, 9).Value = 0.4696472643697726
$ws.Cells.Item(9, 10).Value = 0.4696472643697726
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 4.873409
$ws.Cells.Item(9, 14).Value = 14.620227
$ws.Cells.Item(9, 15).Value = 0.2292080450398417
$ws.Cells.Item(9, 16).Value = 0.2292080450398417
$ws.Cells.Item(9, 17).Value = 37.50484920992599
$ws.Cells.Item(9, 18).Value = 337.543642889334
$ws.Cells.Item(9, 19).Value = 0.1076469313245053
$ws.Cells.Item(9, 20).Value = 0.1076469313245053

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lgi4"
$ws.Cells.Item(10, 3).Value = "Adam22"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 7.695813999999999
$ws.Cells.Item(10, 8).Value = 23.087442
$ws.Cells.Item(10, 9).Value = 0.4696472643697726
$ws.Cells.Item(10, 10).Value = 0.4696472643697726
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 15.55565633333333
$ws.Cells.Item(10, 14).Value = 46.66696899999999
$ws.Cells.Item(10, 15).Value = 0.7316196070296923
$ws.Cells.Item(10, 16).Value = 0.7316196070296922
$ws.Cells.Item(10, 17).Value = 119.7134377892553
$ws.Cells.Item(10, 18).Value = 1077.420940103298
$ws.Cells.Item(10, 19).Value = 0.343603147000783
$ws.Cells.Item(10, 20).Value = 0.343603147000783
